$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 274
$ws.Range("C3").Value = 171363
$ws.Range("C4").Value = 162164
$ws.Range("C5").Value = 9199
$ws.Range("C8").Value = 65.76000000000001
